$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the negative-test data strings for the Accounts module (adds a new
# "Gnukhata_id" negative-accountname case and reworks the negative-groupcode
# case to use "connect_id").
$ws.Range("B33").Value = ', ,%?,Gnukhata_id,@?$1,!@#$%.?,invalid,12str,v2.1,+'
$ws.Range("B34").Value = ', ,%?,Gnukhata_id,@?$1,!@#$%.?,invalid,12str,v2.1,-'
$ws.Range("B35").Value = ', ,testing,%?,@?$1,!@#$%.?,invalid,12str,v2.1,connect_id,%'

# Leave the cursor where the author left it when they committed the change.
[void]$ws.Range("B36").Select()
